$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute("amp; par")
Write-Output ("found: " + $found)

$parRng = $d.Range($rng.End - 4, $rng.End)
Write-Output ("parRng text: [" + $parRng.Text + "]")

$insRng = $parRng.Duplicate
$insRng.Collapse(0)
$insRng.InsertAfter("<lb/>")
$insRng.Font.Name = "Courier New"
$insRng.Font.Color = 11119017
$insRng.Font.Size = 9

Write-Output ("insRng text: [" + $insRng.Text + "]")
Write-Output ("insRng font name: " + $insRng.Font.Name)
Write-Output ("insRng font color: " + $insRng.Font.Color)
Write-Output ("insRng font size: " + $insRng.Font.Size)
